# Updates cryptos list price (D) and 1h volume-change (E) columns.
# D-column values are numeric-looking text (e.g. "1.00", "0.0000280") that must
# stay as literal text, so we force the Text number format before assigning the
# value and then clear the style afterwards so no stray formatting is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.176.19"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.785.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.05%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.785.42"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.03%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -0.55%  "
$ws.Range("E10").Value = "  -3.13%  "
$ws.Range("E11").Value = "  +0.10%  "
$ws.Range("E12").Value = "  -2.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000280"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.62%  "
$ws.Range("E14").Value = "  -1.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.421.76"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.786.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.75"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.160.48"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.02%  "
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "467.65"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.719"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000150"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -8.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.93"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.42"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.31%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.936.08"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.91"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.20%  "
$ws.Range("E32").Value = "  -2.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.59"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.56%  "
$ws.Range("E34").Value = "  -3.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.22"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.742.98"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.75"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.31%  "
$ws.Range("E38").Value = "  -1.14%  "
$ws.Range("E39").Value = "  -0.73%  "
$ws.Range("E40").Value = "  -1.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.87"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.93%  "
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("E43").Value = "  -2.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.68"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.56%  "
$ws.Range("E46").Value = "  -2.64%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "406.75"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.94%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "45.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000273"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -10.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "40.41"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "143.03"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.39%  "
